$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string text updates (header volume number + report week dates) ---
# "Volume 32   Number  10" -> "Volume 32   Number  11"
$volCell = $ws.Range("A8")
$volCell.Characters(21,2).Text = "11"
$volSeg1 = $volCell.Characters(1,20)
$volSeg1.Font.Size = 10
$volSeg1.Font.Name = "Andale WT"
$volSeg2 = $volCell.Characters(21,2)
$volSeg2.Font.Size = 10
$volSeg2.Font.Name = "Andale WT"

# "Report Covering the Week  3/3/2025  Through  3/9/2025"
# -> "Report Covering the Week  3/10/2025  Through  3/16/2025"
$wkCell = $ws.Range("C9")
$wkCell.Characters(46,8).Text = "3/16/2025"
$wkCell.Characters(27,8).Text = "3/10/2025"
$wkSeg1 = $wkCell.Characters(1,26)
$wkSeg1.Font.Size = 10
$wkSeg1.Font.Name = "Andale WT"
$wkSeg2 = $wkCell.Characters(27,9)
$wkSeg2.Font.Size = 10
$wkSeg2.Font.Name = "Andale WT"
$wkSeg3 = $wkCell.Characters(36,11)
$wkSeg3.Font.Size = 10
$wkSeg3.Font.Name = "Andale WT"
$wkSeg4 = $wkCell.Characters(47,9)
$wkSeg4.Font.Size = 10
$wkSeg4.Font.Name = "Andale WT"

# --- Column E width: new data no longer needs the extra bestFit width ---
$ws.Columns.Item(5).ColumnWidth = 5.43

# --- Cells that change "type" (blank-placeholder text <-> real number) ---
# Copy an existing cell with the right style first (so the destination keeps the
# correct number format / font), then overwrite with the new literal value.
$ws.Range("C15").Copy($ws.Range("D15"))
$ws.Range("D15").Value = 1

$ws.Range("H15").Copy($ws.Range("E15"))
$ws.Range("E15").Value = 0

$ws.Range("D20").Copy($ws.Range("C20"))
# "C20" now holds literal text "0" copied from D20

$ws.Range("C27").Copy($ws.Range("D27"))
$ws.Range("D27").Value = 1

$ws.Range("H15").Copy($ws.Range("E27"))
$ws.Range("E27").Value = 0

$ws.Range("C31").Copy($ws.Range("D31"))
# "D31" now holds literal text "0" copied from C31

$ws.Range("E14").Copy($ws.Range("E31"))
# "E31" now holds literal text "***.*" copied from E14

# --- Plain numeric value updates ---
$ws.Range("I15").Value = 6
$ws.Range("J15").Value = 5
$ws.Range("K15").Value = 20
$ws.Range("L15").Value = 200
$ws.Range("N15").Value = 500
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 10
$ws.Range("I16").Value = 23
$ws.Range("J16").Value = 27
$ws.Range("K16").Value = -14.814814814814
$ws.Range("L16").Value = -14.814814814814
$ws.Range("M16").Value = 130
$ws.Range("N16").Value = -88.144329896907
$ws.Range("C17").Value = 2
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 11
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = -21.428571428571
$ws.Range("I17").Value = 22
$ws.Range("J17").Value = 29
$ws.Range("K17").Value = -24.137931034482
$ws.Range("L17").Value = -18.518518518518
$ws.Range("M17").Value = 83.333333333333
$ws.Range("N17").Value = -42.105263157894
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 25
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 24
$ws.Range("H18").Value = -50
$ws.Range("I18").Value = 33
$ws.Range("J18").Value = 54
$ws.Range("K18").Value = -38.888888888888
$ws.Range("L18").Value = 3.125
$ws.Range("M18").Value = -21.428571428571
$ws.Range("N18").Value = -82.065217391304
$ws.Range("C19").Value = 27
$ws.Range("D19").Value = 25
$ws.Range("E19").Value = 8
$ws.Range("F19").Value = 105
$ws.Range("G19").Value = 89
$ws.Range("H19").Value = 17.977528089887
$ws.Range("I19").Value = 257
$ws.Range("J19").Value = 217
$ws.Range("K19").Value = 18.433179723502
$ws.Range("L19").Value = 15.2466367713
$ws.Range("M19").Value = 15.765765765765
$ws.Range("N19").Value = -61.755952380952
$ws.Range("F20").Value = 2
$ws.Range("M20").Value = 66.666666666666
$ws.Range("N20").Value = -97.282608695652
$ws.Range("C21").Value = 36
$ws.Range("D21").Value = 33
$ws.Range("E21").Value = 9.090909090909
$ws.Range("F21").Value = 141
$ws.Range("G21").Value = 138
$ws.Range("H21").Value = 2.173913043478
$ws.Range("I21").Value = 346
$ws.Range("J21").Value = 336
$ws.Range("K21").Value = 2.97619047619
$ws.Range("L21").Value = 8.80503144654
$ws.Range("M21").Value = 19.723183391003
$ws.Range("N21").Value = -72.96875
$ws.Range("C22").Value = 2
$ws.Range("E22").Value = 100
$ws.Range("F22").Value = 11
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = 83.333333333333
$ws.Range("I22").Value = 20
$ws.Range("J22").Value = 15
$ws.Range("K22").Value = 33.333333333333
$ws.Range("L22").Value = -4.761904761904
$ws.Range("M22").Value = 11.111111111111
$ws.Range("C24").Value = 69
$ws.Range("D24").Value = 73
$ws.Range("E24").Value = -5.479452054794
$ws.Range("F24").Value = 273
$ws.Range("G24").Value = 334
$ws.Range("H24").Value = -18.263473053892
$ws.Range("I24").Value = 727
$ws.Range("J24").Value = 838
$ws.Range("K24").Value = -13.245823389021
$ws.Range("L24").Value = -8.091024020227
$ws.Range("M24").Value = 143.959731543624
$ws.Range("C25").Value = 67
$ws.Range("D25").Value = 74
$ws.Range("E25").Value = -9.459459459459
$ws.Range("F25").Value = 278
$ws.Range("G25").Value = 334
$ws.Range("H25").Value = -16.766467065868
$ws.Range("I25").Value = 731
$ws.Range("J25").Value = 840
$ws.Range("K25").Value = -12.97619047619
$ws.Range("L25").Value = -9.417596034696
$ws.Range("C26").Value = 13
$ws.Range("D26").Value = 15
$ws.Range("E26").Value = -13.333333333333
$ws.Range("F26").Value = 33
$ws.Range("H26").Value = -8.333333333333
$ws.Range("I26").Value = 93
$ws.Range("J26").Value = 84
$ws.Range("K26").Value = 10.714285714285
$ws.Range("L26").Value = 47.619047619047
$ws.Range("M26").Value = 78.846153846153
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = 200
$ws.Range("I27").Value = 8
$ws.Range("J27").Value = 6
$ws.Range("K27").Value = 33.333333333333
$ws.Range("L27").Value = 300
$ws.Range("C28").Value = 4
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = 33.333333333333
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 13
$ws.Range("H28").Value = -46.153846153846
$ws.Range("I28").Value = 15
$ws.Range("J28").Value = 22
$ws.Range("K28").Value = -31.818181818181
$ws.Range("L28").Value = -6.25
$ws.Range("F31").Value = 2
$ws.Range("H31").Value = 100
$ws.Range("I31").Value = 4
$ws.Range("K31").Value = 300
$ws.Range("L31").Value = -20
